$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Straightforward text/value updates (Coin, Link, Volume columns, and
#     Price values that are not valid numeric literals) ---
$ws.Cells.Item(2, 4).Value = "30.742.02"
$ws.Cells.Item(2, 5).Value = "  +0.77%  "
$ws.Cells.Item(3, 4).Value = "1.953.42"
$ws.Cells.Item(3, 5).Value = "  +3.44%  "
$ws.Cells.Item(4, 5).Value = "  -0.34%  "
$ws.Cells.Item(5, 5).Value = "  +3.41%  "
$ws.Cells.Item(6, 5).Value = "  +35.06%  "
$ws.Cells.Item(7, 5).Value = "  -0.33%  "
$ws.Cells.Item(8, 5).Value = "  +9.69%  "
$ws.Cells.Item(9, 5).Value = "  +13.00%  "
$ws.Cells.Item(10, 5).Value = "  +5.29%  "
$ws.Cells.Item(11, 5).Value = "  +11.62%  "
$ws.Cells.Item(12, 5).Value = "  +5.58%  "
$ws.Cells.Item(13, 5).Value = "  +2.78%  "
$ws.Cells.Item(14, 4).Value = "1.938.13"
$ws.Cells.Item(14, 5).Value = "  +2.60%  "
$ws.Cells.Item(15, 5).Value = "  +3.09%  "
$ws.Cells.Item(16, 5).Value = "  -1.57%  "
$ws.Cells.Item(17, 4).Value = "30.767.40"
$ws.Cells.Item(17, 5).Value = "  +0.86%  "
$ws.Cells.Item(18, 5).Value = "  +5.89%  "
$ws.Cells.Item(19, 5).Value = "  +2.73%  "
$ws.Cells.Item(20, 5).Value = "  +6.31%  "
$ws.Cells.Item(21, 4).Value = "2.189.51"
$ws.Cells.Item(21, 5).Value = "  +3.17%  "
$ws.Cells.Item(22, 5).Value = "  -0.36%  "
$ws.Cells.Item(23, 5).Value = "  -0.43%  "
$ws.Cells.Item(24, 5).Value = "  +5.55%  "
$ws.Cells.Item(25, 5).Value = "  +4.22%  "
$ws.Cells.Item(26, 5).Value = "  +0.85%  "
$ws.Cells.Item(27, 5).Value = "  +2.64%  "
$ws.Cells.Item(28, 5).Value = "  +11.69%  "
$ws.Cells.Item(29, 5).Value = "  +17.75%  "
$ws.Cells.Item(30, 5).Value = "  +1.43%  "
$ws.Cells.Item(31, 5).Value = "  +5.23%  "
$ws.Cells.Item(32, 5).Value = "  +4.43%  "
$ws.Cells.Item(33, 5).Value = "  +5.17%  "
$ws.Cells.Item(34, 5).Value = "  +2.77%  "
$ws.Cells.Item(35, 5).Value = "  +5.54%  "
$ws.Cells.Item(36, 5).Value = "  +4.23%  "
$ws.Cells.Item(37, 5).Value = "  +0.12%  "
$ws.Cells.Item(38, 5).Value = "  +3.16%  "
$ws.Cells.Item(39, 5).Value = "  +2.68%  "
$ws.Cells.Item(40, 2).Value = "FraxShare"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(40, 5).Value = "  +4.38%  "
$ws.Cells.Item(41, 2).Value = "Aave"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Cells.Item(41, 5).Value = "  +2.12%  "
$ws.Cells.Item(42, 5).Value = "  +7.21%  "
$ws.Cells.Item(43, 5).Value = "  +1.38%  "
$ws.Cells.Item(44, 5).Value = "  +2.22%  "
$ws.Cells.Item(45, 5).Value = "  -0.29%  "
$ws.Cells.Item(46, 5).Value = "  +5.12%  "
$ws.Cells.Item(47, 5).Value = "  +1.01%  "
$ws.Cells.Item(48, 5).Value = "  +4.32%  "
$ws.Cells.Item(49, 5).Value = "  +2.26%  "
$ws.Cells.Item(50, 5).Value = "  +5.76%  "
$ws.Cells.Item(51, 5).Value = "  +0.31%  "

# --- Price column updates whose new text looks like a plain number. Excel
#     would silently convert a direct .Value assignment into a numeric cell,
#     dropping significant trailing/leading zeros and changing the stored
#     type away from text. Route these through a scratch cell that is force-
#     typed as text (leading apostrophe) and Copy / PasteSpecial *values only*
#     into the destination so the destination cell keeps its original (default)
#     style while the content lands as text, matching the source data.
$scratch = $ws.Cells.Item(200, 200)
$scratch.Value = "'0.9976"
$scratch.Copy()
$ws.Cells.Item(4, 4).PasteSpecial(-4163)
$scratch.Value = "'252.22"
$scratch.Copy()
$ws.Cells.Item(5, 4).PasteSpecial(-4163)
$scratch.Value = "'0.6337"
$scratch.Copy()
$ws.Cells.Item(6, 4).PasteSpecial(-4163)
$scratch.Value = "'0.9973"
$scratch.Copy()
$ws.Cells.Item(7, 4).PasteSpecial(-4163)
$scratch.Value = "'0.3181"
$scratch.Copy()
$ws.Cells.Item(8, 4).PasteSpecial(-4163)
$scratch.Value = "'25.23"
$scratch.Copy()
$ws.Cells.Item(9, 4).PasteSpecial(-4163)
$scratch.Value = "'0.06842"
$scratch.Copy()
$ws.Cells.Item(10, 4).PasteSpecial(-4163)
$scratch.Value = "'0.8131"
$scratch.Copy()
$ws.Cells.Item(11, 4).PasteSpecial(-4163)
$scratch.Value = "'101.13"
$scratch.Copy()
$ws.Cells.Item(12, 4).PasteSpecial(-4163)
$scratch.Value = "'0.07965"
$scratch.Copy()
$ws.Cells.Item(13, 4).PasteSpecial(-4163)
$scratch.Value = "'278.00"
$scratch.Copy()
$ws.Cells.Item(16, 4).PasteSpecial(-4163)
$scratch.Value = "'13.82"
$scratch.Copy()
$ws.Cells.Item(18, 4).PasteSpecial(-4163)
$scratch.Value = "'0.000007695"
$scratch.Copy()
$ws.Cells.Item(19, 4).PasteSpecial(-4163)
$scratch.Value = "'5.609"
$scratch.Copy()
$ws.Cells.Item(20, 4).PasteSpecial(-4163)
$scratch.Value = "'0.9972"
$scratch.Copy()
$ws.Cells.Item(22, 4).PasteSpecial(-4163)
$scratch.Value = "'0.9967"
$scratch.Copy()
$ws.Cells.Item(23, 4).PasteSpecial(-4163)
$scratch.Value = "'6.625"
$scratch.Copy()
$ws.Cells.Item(24, 4).PasteSpecial(-4163)
$scratch.Value = "'9.480"
$scratch.Copy()
$ws.Cells.Item(25, 4).PasteSpecial(-4163)
$scratch.Value = "'165.00"
$scratch.Copy()
$ws.Cells.Item(26, 4).PasteSpecial(-4163)
$scratch.Value = "'19.47"
$scratch.Copy()
$ws.Cells.Item(27, 4).PasteSpecial(-4163)
$scratch.Value = "'2.119"
$scratch.Copy()
$ws.Cells.Item(28, 4).PasteSpecial(-4163)
$scratch.Value = "'0.1145"
$scratch.Copy()
$ws.Cells.Item(29, 4).PasteSpecial(-4163)
$scratch.Value = "'1.354"
$scratch.Copy()
$ws.Cells.Item(30, 4).PasteSpecial(-4163)
$scratch.Value = "'1.548"
$scratch.Copy()
$ws.Cells.Item(31, 4).PasteSpecial(-4163)
$scratch.Value = "'4.480"
$scratch.Copy()
$ws.Cells.Item(32, 4).PasteSpecial(-4163)
$scratch.Value = "'4.338"
$scratch.Copy()
$ws.Cells.Item(33, 4).PasteSpecial(-4163)
$scratch.Value = "'1.190"
$scratch.Copy()
$ws.Cells.Item(35, 4).PasteSpecial(-4163)
$scratch.Value = "'0.7241"
$scratch.Copy()
$ws.Cells.Item(36, 4).PasteSpecial(-4163)
$scratch.Value = "'2.719"
$scratch.Copy()
$ws.Cells.Item(37, 4).PasteSpecial(-4163)
$scratch.Value = "'0.01958"
$scratch.Copy()
$ws.Cells.Item(38, 4).PasteSpecial(-4163)
$scratch.Value = "'2.921"
$scratch.Copy()
$ws.Cells.Item(39, 4).PasteSpecial(-4163)
$scratch.Value = "'6.483"
$scratch.Copy()
$ws.Cells.Item(40, 4).PasteSpecial(-4163)
$scratch.Value = "'77.50"
$scratch.Copy()
$ws.Cells.Item(41, 4).PasteSpecial(-4163)
$scratch.Value = "'0.4567"
$scratch.Copy()
$ws.Cells.Item(42, 4).PasteSpecial(-4163)
$scratch.Value = "'2.030"
$scratch.Copy()
$ws.Cells.Item(43, 4).PasteSpecial(-4163)
$scratch.Value = "'0.8449"
$scratch.Copy()
$ws.Cells.Item(44, 4).PasteSpecial(-4163)
$scratch.Value = "'0.9973"
$scratch.Copy()
$ws.Cells.Item(45, 4).PasteSpecial(-4163)
$scratch.Value = "'10.03"
$scratch.Copy()
$ws.Cells.Item(46, 4).PasteSpecial(-4163)
$scratch.Value = "'102.49"
$scratch.Copy()
$ws.Cells.Item(47, 4).PasteSpecial(-4163)
$scratch.Value = "'7.288"
$scratch.Copy()
$ws.Cells.Item(48, 4).PasteSpecial(-4163)
$scratch.Value = "'35.96"
$scratch.Copy()
$ws.Cells.Item(49, 4).PasteSpecial(-4163)
$scratch.Value = "'0.4170"
$scratch.Copy()
$ws.Cells.Item(50, 4).PasteSpecial(-4163)
$scratch.Value = "'920.12"
$scratch.Copy()
$ws.Cells.Item(51, 4).PasteSpecial(-4163)
$scratch.Clear()
$excel.CutCopyMode = 0

